$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing header "Category" in A1, using the same formatting as the
# rest of the header row (e.g. B1), which carries the bold/centered/bordered style.
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats

# The category cells A2:A46 previously (incorrectly) carried that same header
# style; reset them to the plain/default formatting used by the rest of the
# data cells (e.g. B2), leaving their text values untouched.
$ws.Range("B2").Copy()
$ws.Range("A2:A46").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
